$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain numeric-looking string need
# to be forced to Text format first, otherwise Excel would silently turn them
# into real numbers instead of keeping them as text like the original data.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "30.125.46"
$ws.Range("E2").Value = "  +10.04%  "

$ws.Range("D3").Value = "1.869.77"
$ws.Range("E3").Value = "  +6.89%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.46%  "

$ws.Range("D5").Value = "250.09"
$ws.Range("E5").Value = "  +3.47%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").Value = "0.4964"
$ws.Range("E7").Value = "  +3.08%  "

$ws.Range("D8").Value = "45.05"
$ws.Range("E8").Value = "  +8.94%  "

$ws.Range("D9").Value = "0.2846"
$ws.Range("E9").Value = "  +9.08%  "

$ws.Range("D10").Value = "0.06542"
$ws.Range("E10").Value = "  +6.11%  "

$ws.Range("D11").Value = "1.869.24"
$ws.Range("E11").Value = "  +7.03%  "

$ws.Range("D12").Value = "16.99"
$ws.Range("E12").Value = "  +5.34%  "

$ws.Range("D13").Value = "0.07195"
$ws.Range("E13").Value = "  +3.68%  "

$ws.Range("D14").Value = "0.6635"
$ws.Range("E14").Value = "  +9.70%  "

$ws.Range("D15").Value = "85.40"
$ws.Range("E15").Value = "  +10.64%  "

$ws.Range("D16").Value = "4.808"
$ws.Range("E16").Value = "  +7.66%  "

$ws.Range("D17").Value = "30.150.02"
$ws.Range("E17").Value = "  +10.25%  "

$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "0.000007532"
$ws.Range("E19").Value = "  +6.11%  "

$ws.Range("D20").Value = "12.64"
$ws.Range("E20").Value = "  +10.31%  "

$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "2.105.52"

$ws.Range("D23").Value = "4.714"
$ws.Range("E23").Value = "  +6.26%  "

$ws.Range("D24").Value = "5.500"
$ws.Range("E24").Value = "  +7.42%  "

$ws.Range("D25").Value = "8.985"
$ws.Range("E25").Value = "  +6.51%  "

$ws.Range("D26").Value = "144.17"
$ws.Range("E26").Value = "  +1.23%  "

$ws.Range("D27").Value = "133.72"
$ws.Range("E27").Value = "  +24.39%  "

$ws.Range("D28").Value = "16.72"
$ws.Range("E28").Value = "  +9.52%  "

$ws.Range("D29").Value = "1.938"
$ws.Range("E29").Value = "  +5.68%  "

$ws.Range("D30").Value = "1.400"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").Value = "4.247"
$ws.Range("E31").Value = "  +7.58%  "

$ws.Range("D32").Value = "0.08584"
$ws.Range("E32").Value = "  +7.71%  "

$ws.Range("D33").Value = "3.872"
$ws.Range("E33").Value = "  +5.53%  "

$ws.Range("D34").Value = "0.05050"
$ws.Range("E34").Value = "  +7.07%  "

$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  +11.12%  "

$ws.Range("D36").Value = "0.6835"
$ws.Range("E36").Value = "  +10.42%  "

$ws.Range("D37").Value = "2.688"
$ws.Range("E37").Value = "  +2.69%  "

$ws.Range("D38").Value = "2.320"
$ws.Range("E38").Value = "  +14.82%  "

$ws.Range("D39").Value = "2.746"
$ws.Range("E39").Value = "  +7.80%  "

$ws.Range("D40").Value = "0.9546"
$ws.Range("E40").Value = "  +2.80%  "

$ws.Range("D41").Value = "0.01631"
$ws.Range("E41").Value = "  +9.08%  "

$ws.Range("D42").Value = "6.157"
$ws.Range("E42").Value = "  +7.53%  "

$ws.Range("D43").Value = "0.9997"
$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "103.18"
$ws.Range("E44").Value = "  +3.79%  "

$ws.Range("D45").Value = "0.4173"
$ws.Range("E45").Value = "  +8.59%  "

$ws.Range("D46").Value = "7.407"
$ws.Range("E46").Value = "  +7.60%  "

$ws.Range("D47").Value = "0.1249"
$ws.Range("E47").Value = "  +8.18%  "

$ws.Range("D48").Value = "0.05628"
$ws.Range("E48").Value = "  +5.02%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "32.43"
$ws.Range("E49").Value = "  +8.77%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.297"
$ws.Range("E50").Value = "  +6.16%  "

$ws.Range("E51").Value = "  +6.90%  "
